$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.850.20'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '1.803.76'
$ws.Range("E3").Value = '  -0.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.60%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.10'
$ws.Range("E5").Value = '  +0.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.56%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4593'
$ws.Range("E7").Value = '  +4.24%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3741'
$ws.Range("E8").Value = '  +0.81%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07364'
$ws.Range("E9").Value = '  -0.91%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8603'
$ws.Range("E10").Value = '  -0.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.43'
$ws.Range("E11").Value = '  -1.23%  '

$ws.Range("D12").Value = '1.815.97'
$ws.Range("E12").Value = '  +0.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.626'
$ws.Range("E13").Value = '  -0.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.363'
$ws.Range("E14").Value = '  +1.52%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07086'
$ws.Range("E15").Value = '  +0.23%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.37'
$ws.Range("E16").Value = '  -1.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("E17").Value = '  +0.62%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008687'
$ws.Range("E18").Value = '  -0.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  +0.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.82'
$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("D21").Value = '26.885.47'
$ws.Range("E21").Value = '  -0.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.301'
$ws.Range("E22").Value = '  +2.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.81'
$ws.Range("E23").Value = '  -0.23%  '

$ws.Range("D24").Value = '2.035.75'
$ws.Range("E24").Value = '  +0.75%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.907'
$ws.Range("E25").Value = '  -3.74%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.36'
$ws.Range("E26").Value = '  +0.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.194'
$ws.Range("E27").Value = '  -0.65%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.50'
$ws.Range("E28").Value = '  +0.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.238'
$ws.Range("E29").Value = '  +0.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.30'
$ws.Range("E30").Value = '  -1.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08881'
$ws.Range("E31").Value = '  +1.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7687'
$ws.Range("E32").Value = '  +2.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.163'
$ws.Range("E33").Value = '  -0.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.489'
$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.897'
$ws.Range("E35").Value = '  +0.31%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.005'
$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.125'
$ws.Range("E37").Value = '  +2.51%  '

$ws.Range("E38").Value = '  -0.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05218'
$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.907'
$ws.Range("E40").Value = '  +3.09%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.203'
$ws.Range("E41").Value = '  +1.74%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5259'
$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.333'
$ws.Range("E43").Value = '  +8.82%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1667'
$ws.Range("E44").Value = '  -1.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.576'
$ws.Range("E45").Value = '  +0.60%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4990'
$ws.Range("E46").Value = '  +0.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.24'
$ws.Range("E47").Value = '  -1.98%  '

$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.005'
$ws.Range("E48").Value = '  +0.59%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '104.22'
$ws.Range("E49").Value = '  -0.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.662'
$ws.Range("E50").Value = '  -0.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06315'
$ws.Range("E51").Value = '  -0.48%  '
